$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New price-report entry inserted as row 6 ---
# Shifts old rows 6-11 down to 7-12.
$ws.Rows(6).Insert()

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44859
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100104
$ws.Range("H6").Value = "Frutos de pepita"
$ws.Range("I6").Value = 100104004
$ws.Range("J6").Value = "Níspero"
$ws.Range("K6").Value = "Californiana(o)"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("Q6").Value = "$/bandeja 5 kilos"
$ws.Range("R6").Value = "Provincia de Quillota"
$ws.Range("S6").Value = 4000
$ws.Range("T6").Value = 5

# --- Second new price-report entry inserted as row 12 ---
# Shifts the (now) row 12 (old row 11) down to row 13.
$ws.Rows(12).Insert()

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 44858
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100104
$ws.Range("H12").Value = "Frutos de pepita"
$ws.Range("I12").Value = 100104004
$ws.Range("J12").Value = "Níspero"
$ws.Range("K12").Value = "Californiana(o)"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 90
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("Q12").Value = "$/bandeja 5 kilos"
$ws.Range("R12").Value = "Provincia de Quillota"
$ws.Range("S12").Value = 4000
$ws.Range("T12").Value = 5
